$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '46.047.21'
$ws.Range("E2").Value = '  -1.63%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.371.12'
$ws.Range("E3").Value = '  +2.55%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '301.62'
$ws.Range("E5").Value = '  +0.61%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '99.03'
$ws.Range("E6").Value = '  -3.49%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.570'
$ws.Range("E7").Value = '  -0.72%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.516'
$ws.Range("E9").Value = '  -2.68%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.53'
$ws.Range("E10").Value = '  -6.71%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0796'
$ws.Range("E11").Value = '  -1.16%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.14'
$ws.Range("E12").Value = '  -3.60%  '
$ws.Range("E13").Value = '  -0.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.733.79'
$ws.Range("E14").Value = '  +2.76%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.355.01'
$ws.Range("E15").Value = '  +2.11%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.816'
$ws.Range("E16").Value = '  -0.96%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.64'
$ws.Range("E17").Value = '  -2.65%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '45.984.96'
$ws.Range("E18").Value = '  -1.65%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.89'
$ws.Range("E19").Value = '  -3.36%  '
$ws.Range("E20").Value = '  +1.59%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.05'
$ws.Range("E21").Value = '  -1.87%  '
$ws.Range("E22").Value = '  +0.44%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '245.25'
$ws.Range("E23").Value = '  -1.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.82'
$ws.Range("E24").Value = '  -4.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.93'
$ws.Range("E26").Value = '  -2.82%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '39.85'
$ws.Range("E27").Value = '  -7.90%  '
$ws.Range("E28").Value = '  -3.27%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.81'
$ws.Range("E29").Value = '  -1.37%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.79'
$ws.Range("E30").Value = '  +21.54%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '21.02'
$ws.Range("E31").Value = '  +3.90%  '
$ws.Range("E32").Value = '  +6.60%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.53'
$ws.Range("E33").Value = '  -4.35%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '146.57'
$ws.Range("E34").Value = '  -0.37%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0778'
$ws.Range("E35").Value = '  -2.99%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.113'
$ws.Range("E36").Value = '  +0.78%  '
$ws.Range("E37").Value = '  +5.16%  '
$ws.Range("E38").Value = '  -2.94%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.02'
$ws.Range("E39").Value = '  -5.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.94'
$ws.Range("E40").Value = '  -2.42%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0301'
$ws.Range("E41").Value = '  -2.18%  '
$ws.Range("E42").Value = '  -6.63%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.920.05'
$ws.Range("E43").Value = '  +4.27%  '
$ws.Range("E44").Value = '  +0.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.02'
$ws.Range("E45").Value = '  +2.67%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.79'
$ws.Range("E46").Value = '  -11.16%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.39'
$ws.Range("E47").Value = '  +4.81%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.187'
$ws.Range("E48").Value = '  -5.55%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '98.11'
$ws.Range("E49").Value = '  +0.38%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.604.73'
$ws.Range("E50").Value = '  +2.61%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '69.06'
$ws.Range("E51").Value = '  -8.51%  '
